$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up formatting first (copy formats from cells that already carry
# the target look), then overwrite the cell text. ---

# I3 needs to drop back to the default/unstyled look (like Q10, which is
# a plain text cell with no explicit style).
$ws.Range("Q10").Copy()
$ws.Range("I3").PasteSpecial(-4122)

# J3 needs the "INT input" look (style used at O11).
$ws.Range("O11").Copy()
$ws.Range("J3").PasteSpecial(-4122)

# I4 / J4 need the "GND" look (style used at H3/H4).
$ws.Range("H3").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("H4").Copy()
$ws.Range("J4").PasteSpecial(-4122)

# I5 / J5 need the generic analog-signal look (style used at G3/L3).
$ws.Range("G3").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("L3").Copy()
$ws.Range("J5").PasteSpecial(-4122)

# I6 / J6 need the "3.3V" look (style used at F3/F4).
$ws.Range("F3").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("J6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Now set the new pin-label text for connector 1, pins 1-4 ---

# Row 3 (pin 1)
$ws.Range("G3").Value = "VTemp"
$ws.Range("I3").Value = "Ign Cut"
$ws.Range("J3").Value = "INT input"
$ws.Range("L3").Value = "Travel 4"

# Row 4 (pin 2)
$ws.Range("I4").Value = "GND"
$ws.Range("J4").Value = "GND"
$ws.Range("L4").Value = "Travel 3"

# Row 5 (pin 3)
$ws.Range("I5").Value = "BrakeP"
$ws.Range("J5").Value = "SpeedPos"
$ws.Range("L5").Value = "Travel 1"

# Row 6 (pin 4)
$ws.Range("I6").Value = "3.3V"
$ws.Range("J6").Value = "3.3V"
$ws.Range("L6").Value = "Travel 2"

# --- Update the recorded selection on the sheet view ---
$ws.Range("N5").Select()
